$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.655.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.01%  "

$ws.Range("D3").Value = "'2.624.53"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.00%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'595.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.33%  "

$ws.Range("D6").Value = "'150.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.09%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "'0.588"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.15%  "

$ws.Range("E9").Value = "  -0.09%  "

$ws.Range("E10").Value = "  +1.90%  "

$ws.Range("E11").Value = "  +2.82%  "

$ws.Range("E12").Value = "  -1.15%  "

$ws.Range("D13").Value = "'27.67"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.10%  "

$ws.Range("D14").Value = "'3.098.39"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.87%  "

$ws.Range("D15").Value = "'63.473.24"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.03%  "

$ws.Range("E16").Value = "  +2.05%  "

$ws.Range("D17").Value = "'2.638.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.64%  "

$ws.Range("D18").Value = "'12.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +7.02%  "

$ws.Range("E19").Value = "  +2.05%  "

$ws.Range("D20").Value = "'347.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.15%  "

$ws.Range("D21").Value = "'6.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.22%  "

$ws.Range("E22").Value = "  -0.14%  "

$ws.Range("D23").Value = "'5.73"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.50%  "

$ws.Range("D24").Value = "'66.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.72%  "

$ws.Range("E25").Value = "  +11.19%  "

$ws.Range("D26").Value = "'1.68"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.44%  "

$ws.Range("D27").Value = "'9.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.41%  "

$ws.Range("D28").Value = "'576.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.31%  "

$ws.Range("D29").Value = "'8.22"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.61%  "

$ws.Range("E30").Value = "  +0.06%  "

$ws.Range("E31").Value = "  +0.11%  "

$ws.Range("D32").Value = "'2.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.05%  "

$ws.Range("D33").Value = "'0.0₃0845"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.38%  "

$ws.Range("E34").Value = "  -0.54%  "

$ws.Range("D35").Value = "'5.25"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.47%  "

$ws.Range("D36").Value = "'168.75"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.03%  "

$ws.Range("E37").Value = "  +0.24%  "

$ws.Range("D38").Value = "'1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.02%  "

$ws.Range("E39").Value = "  +0.24%  "

$ws.Range("D40").Value = "'19.35"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.14%  "

$ws.Range("D42").Value = "'168.82"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.26%  "

$ws.Range("D43").Value = "'39.90"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.42%  "

$ws.Range("E44").Value = "  +3.89%  "

$ws.Range("D45").Value = "'0.0601"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.96%  "

$ws.Range("D46").Value = "'21.36"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.92%  "

$ws.Range("E47").Value = "  -0.79%  "

$ws.Range("E48").Value = "  +0.87%  "

$ws.Range("E49").Value = "  +4.75%  "

$ws.Range("D50").Value = "'0.0964"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.13%  "

$ws.Range("E51").Value = "  +2.46%  "

